# The deck currently has two DrawingML theme parts:
#   ppt/theme/theme1.xml -> bound to the slide master ("Integral" / "Red Violet" colours)
#   ppt/theme/theme2.xml -> bound to the notes master  ("Office Theme" / "Office" colours)
#
# The authored edit swaps the contents of those two theme parts: the slide
# master's theme becomes the stock "Office" colour scheme, and the notes
# master's theme becomes the old "Integral" colour scheme.
#
# Apply the colour-scheme half of that swap through the PowerPoint object
# model: walk the 12 theme colour slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) on the slide master's Theme and set each one to the RGB
# value that the "Office" theme used to have.

function Convert-HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target colour scheme (was theme2.xml's "Office" scheme, now theme1.xml's).
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = Convert-HexToRGB $officeColors[$i - 1]
}
